$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'72.117.69"
$ws.Range("E2").Value = "  +5.27%  "

# Row 3
$ws.Range("D3").Value = "'4.051.49"
$ws.Range("E3").Value = "  +5.10%  "

# Row 5
$ws.Range("D5").Value = "'532.76"
$ws.Range("E5").Value = "  +2.74%  "

# Row 6
$ws.Range("D6").Value = "'152.45"
$ws.Range("E6").Value = "  +8.55%  "

# Row 7
$ws.Range("D7").Value = "'0.705"
$ws.Range("E7").Value = "  +15.96%  "

# Row 8
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
$ws.Range("D9").Value = "'0.768"
$ws.Range("E9").Value = "  +8.35%  "

# Row 10
$ws.Range("E10").Value = "  +6.32%  "

# Row 11
$ws.Range("D11").Value = "'0.0000334"
$ws.Range("E11").Value = "  +4.47%  "

# Row 12
$ws.Range("D12").Value = "'50.14"
$ws.Range("E12").Value = "  +20.77%  "

# Row 13
$ws.Range("D13").Value = "'11.09"
$ws.Range("E13").Value = "  +8.41%  "

# Row 14
$ws.Range("D14").Value = "'4.704.71"
$ws.Range("E14").Value = "  +5.61%  "

# Row 15
$ws.Range("D15").Value = "'4.077.68"
$ws.Range("E15").Value = "  +6.12%  "

# Row 16
$ws.Range("E16").Value = "  +2.70%  "

# Row 17
$ws.Range("D17").Value = "'21.14"
$ws.Range("E17").Value = "  -2.66%  "

# Row 18
$ws.Range("D18").Value = "'1.23"
$ws.Range("E18").Value = "  +2.42%  "

# Row 19
$ws.Range("D19").Value = "'0.134"
$ws.Range("E19").Value = "  +0.16%  "

# Row 20
$ws.Range("D20").Value = "'72.097.77"
$ws.Range("E20").Value = "  +5.34%  "

# Row 21
$ws.Range("D21").Value = "'439.41"
$ws.Range("E21").Value = "  +5.41%  "

# Row 22
$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").Value = "'3.72"
$ws.Range("E22").Value = "  +8.28%  "

# Row 23
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'100.60"
$ws.Range("E23").Value = "  +16.15%  "

# Row 24
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'4.30"
$ws.Range("E24").Value = "  +8.76%  "

# Row 25
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "'14.86"
$ws.Range("E25").Value = "  +6.42%  "

# Row 26
$ws.Range("D26").Value = "'11.45"
$ws.Range("E26").Value = "  -0.37%  "

# Row 27
$ws.Range("D27").Value = "'11.02"
$ws.Range("E27").Value = "  +4.59%  "

# Row 28
$ws.Range("D28").Value = "'37.50"
$ws.Range("E28").Value = "  +6.16%  "

# Row 29
$ws.Range("D29").Value = "'5.85"
$ws.Range("E29").Value = "  +3.25%  "

# Row 30
$ws.Range("D30").Value = "'3.49"
$ws.Range("E30").Value = "  +26.13%  "

# Row 31
$ws.Range("D31").Value = "'13.76"
$ws.Range("E31").Value = "  +4.82%  "

# Row 32
$ws.Range("E32").Value = "  +6.99%  "

# Row 33
$ws.Range("D33").Value = "'675.35"
$ws.Range("E33").Value = "  +0.52%  "

# Row 34
$ws.Range("D34").Value = "'6.77"
$ws.Range("E34").Value = "  +5.12%  "

# Row 35
$ws.Range("D35").Value = "'66.86"
$ws.Range("E35").Value = "  +2.90%  "

# Row 36
$ws.Range("D36").Value = "'42.91"
$ws.Range("E36").Value = "  +8.26%  "

# Row 37
$ws.Range("D37").Value = "'0.438"
$ws.Range("E37").Value = "  -1.79%  "

# Row 38
$ws.Range("E38").Value = "  +4.60%  "

# Row 39
$ws.Range("E39").Value = "  +7.65%  "

# Row 40
$ws.Range("D40").Value = "'3.44"
$ws.Range("E40").Value = "  -0.19%  "

# Row 41
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.07%  "

# Row 42
$ws.Range("D42").Value = "'0.0503"
$ws.Range("E42").Value = "  +6.26%  "

# Row 43
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.15%  "

# Row 44
$ws.Range("D44").Value = "'3.19"
$ws.Range("E44").Value = "  +2.99%  "

# Row 45
$ws.Range("D45").Value = "'0.154"
$ws.Range("E45").Value = "  +10.88%  "

# Row 46
$ws.Range("D46").Value = "'2.76"
$ws.Range("E46").Value = "  -0.31%  "

# Row 47
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").Value = "'9.58"
$ws.Range("E47").Value = "  +12.82%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'3.38"
$ws.Range("E48").Value = "  +0.27%  "

# Row 49
$ws.Range("E49").Value = "  +5.09%  "

# Row 50
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").Value = "'0.000278"
$ws.Range("E50").Value = "  +1.18%  "

# Row 51
$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").Value = "'3.39"
$ws.Range("E51").Value = "  +4.07%  "
